$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Max Ignacio Lastra Yañez"
$ws.Range("B4").Value = "max.lastray@gmail.com"
$ws.Range("C4").Value = "Atletismo"
$ws.Range("D4").Formula = '="2024-11-03"'
$ws.Range("E4").Value = "15:00"

$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
